$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change ("Generate Report for Handback") drops the handback row for the
# file dcdaf9f5-af2c-4b5f-83bf-f8f5a6c29cdf (row 3 on every sheet) and
# refreshes the handoff/handback timestamps recorded for the remaining
# ab3f5bcd-c748-4009-8642-b2e8ef5797d8 row on the zh-cn and de-de sheets.
# ---------------------------------------------------------------------------

# ---- Overview sheet --------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()
$wsOverview.Rows.Item(3).Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/598af6488b6eca5870a0af7602a5035e891ee0b4/e2e/ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md", "", "", "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md")

# ---- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Rows.Item(3).Delete()

$wsZhCn.Range("E2").Value = "2016-03-19 00:37:35"
$wsZhCn.Range("H2").Value = "2016-03-19 00:37:53"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/598af6488b6eca5870a0af7602a5035e891ee0b4/e2e/ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md", "", "", "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/598af6488b6eca5870a0af7602a5035e891ee0b4/e2e/ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b307fbe86a05fe985feb3b2e43e09f718c632f1e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.zh-cn.xlf", "", "", "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b2f84f44030c36ac45b5a93760524914046f46d4/e2e/ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md", "", "", "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f2460b8cdb19560b145887271dd805820db6991d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.zh-cn.xlf", "", "", "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.zh-cn.xlf")

# ---- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Rows.Item(3).Delete()

$wsDeDe.Range("E2").Value = "2016-03-19 00:37:39"
$wsDeDe.Range("H2").Value = "2016-03-19 00:37:58"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/598af6488b6eca5870a0af7602a5035e891ee0b4/e2e/ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md", "", "", "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/598af6488b6eca5870a0af7602a5035e891ee0b4/e2e/ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f652705f1e3ebd943ee08f144cfdbf903e8884b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.de-de.xlf", "", "", "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/04287f3319d7db557c821b517a75fb89feebb9e6/e2e/ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md", "", "", "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/78a76bb3f1dc7e55788dd8703fce71cf00154b02/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.de-de.xlf", "", "", "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.de-de.xlf")
